$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing merges so the structural edits below are unambiguous.
$ws.Range("A1:I10").UnMerge()

# Drop the now-unused H/I columns entirely.
$ws.Range("H1:I10").Delete()

# Drop rows 5-10; remaining data becomes rows 1-4.
$ws.Range("A5:A10").EntireRow.Delete()

# --- Row 2 (Conduit 1, pull 1) ---
$ws.Range("A2").Value = "Conduit 1"
$ws.Range("B2").Value = "543+00"
$ws.Range("C2").Value = "553+00"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "7C#14"
$ws.Range("F2").Value = "EXPRESS"

# --- Row 3 (Conduit 2, pull 2) ---
$ws.Range("A3").Value = "Conduit 2"
$ws.Range("B3").Value = "543+00"
$ws.Range("C3").Value = "553+00"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "7C#14"
$ws.Range("F3").Value = "LOCAL"

# --- Row 4 (pull 3, merged continuation of Conduit 2) ---
$ws.Range("A4").Value = $null
$ws.Range("B4").Value = $null
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "7C#14"
$ws.Range("F4").Value = $null

# Percentage-looking text ("10.91%", "21.82%") must stay literal text, not be
# reinterpreted as a numeric percent value by Excel's entry auto-detection.
# Write as a formula producing the literal string, then freeze it to a
# plain value via copy / paste-values so no formula or new number format
# sticks around.
$ws.Range("G2").Formula = '="10.91%"'
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0

$ws.Range("G3").Formula = '="21.82%"'
$ws.Range("G3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0

$ws.Range("G4").Value = $null

# Re-merge the Conduit-2 block (rows 3-4) across A, B, C, F, G.
$ws.Range("A3:A4").Merge()
$ws.Range("B3:B4").Merge()
$ws.Range("C3:C4").Merge()
$ws.Range("F3:F4").Merge()
$ws.Range("G3:G4").Merge()

Write-Output $ws.UsedRange.Address()
